$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
# ---- populate 15 data rows (rows 2..16) ----
$ws1.Range("B2:B16").NumberFormat = "@"
$ws1.Range("E2:E16").NumberFormat = "@"
$data = New-Object 'object[,]' 15,8
$data[0,0] = '2024-05-18'
$data[0,1] = '合肥·WA二次元饭局（取消）'
$data[0,2] = '临泉路胜利路交叉路（中环国际大厦对面） 太太满庭芳(胜利路店)'
$data[0,3] = '2024.05.18 14:50-05.18 20:00'
$data[0,4] = 67
$data[0,5] = '不可售'
$data[0,6] = 'https://show.bilibili.com/platform/detail.html?id=83978'
$data[0,7] = '//i1.hdslb.com/bfs/openplatform/202404/wK9Yq9Ta1712657384067.jpeg'
$data[1,0] = '2024-05-18'
$data[1,1] = '合肥·梦时空SPO1动漫展（取消）'
$data[1,2] = '阜阳路16号 银瑞林国际大酒店'
$data[1,3] = '2024.05.18 10:00-05.18 17:00'
$data[1,4] = 126
$data[1,5] = '不可售'
$data[1,6] = 'https://show.bilibili.com/platform/detail.html?id=80207'
$data[1,7] = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'
$data[2,0] = '2024-05-18'
$data[2,1] = '合肥·首届偶活企划——偶像计划-闪耀舞台'
$data[2,2] = '阜阳路16号 银瑞林国际大酒店'
$data[2,3] = '2024.05.18 09:00-05.18 17:00'
$data[2,4] = 75
$data[2,5] = 58
$data[2,6] = 'https://show.bilibili.com/platform/detail.html?id=83891'
$data[2,7] = '//i2.hdslb.com/bfs/openplatform/202404/lfqv8l9Q1712453982625.jpeg'
$data[3,0] = '2024-06-01'
$data[3,1] = '合肥·TH元气动漫游戏嘉年华'
$data[3,2] = '北一环路与胜利路交口西北侧中星城2号楼(地铁1号线长淮站D出口） 曙光薇酒店(合肥站店)'
$data[3,3] = '2024.06.01 10:00-06.01 17:00'
$data[3,4] = 3
$data[3,5] = 50
$data[3,6] = 'https://show.bilibili.com/platform/detail.html?id=85181'
$data[3,7] = '//i0.hdslb.com/bfs/openplatform/202405/QcP0IlNZ1715064886624.jpeg'
$data[4,0] = '2024-06-01'
$data[4,1] = '合肥·运动番only·群青日和'
$data[4,2] = '金寨路287号 合肥明星运动公园'
$data[4,3] = '2024.06.01 09:30-06.01 17:30'
$data[4,4] = 539
$data[4,5] = 80
$data[4,6] = 'https://show.bilibili.com/platform/detail.html?id=83058'
$data[4,7] = '//i2.hdslb.com/bfs/openplatform/202404/Jzeq47lD1714026878824.jpeg'
$data[5,0] = '2024-06-08'
$data[5,1] = '合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~'
$data[5,2] = '锦绣大道3899号 合肥滨湖会展中心'
$data[5,3] = '2024.06.08 09:30-06.09 17:00'
$data[5,4] = 7523
$data[5,5] = 65
$data[5,6] = 'https://show.bilibili.com/platform/detail.html?id=83518'
$data[5,7] = '//i1.hdslb.com/bfs/openplatform/202403/1Sqp42gM1711691520194.jpeg'
$data[6,0] = '2024-06-09'
$data[6,1] = '合肥·第二届华盟动漫次元嘉年华'
$data[6,2] = '常青街道十五里河村合柴1972院内 合肥当代美术馆'
$data[6,3] = '2024.06.09 10:00-06.10 17:00'
$data[6,4] = 474
$data[6,5] = 58
$data[6,6] = 'https://show.bilibili.com/platform/detail.html?id=84081'
$data[6,7] = '//i1.hdslb.com/bfs/openplatform/202404/O5LyHE7j1712732240786.jpeg'
$data[7,0] = '2024-06-09'
$data[7,1] = '合肥·第六届环形宇宙动漫游戏嘉年华内场票·赵成晨'
$data[7,2] = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$data[7,3] = '2024.06.09 09:30-06.09 17:00'
$data[7,4] = 194
$data[7,5] = 238
$data[7,6] = 'https://show.bilibili.com/platform/detail.html?id=84863'
$data[7,7] = '//i1.hdslb.com/bfs/openplatform/202404/I5S4Ih2M1714031127805.jpeg'
$data[8,0] = '2024-06-22'
$data[8,1] = '合肥·Look Look动漫嘉年华'
$data[8,2] = '新站区东方大道288号 少荃体育中心'
$data[8,3] = '2024.06.22 10:00-06.22 17:30'
$data[8,4] = 1069
$data[8,5] = 58
$data[8,6] = 'https://show.bilibili.com/platform/detail.html?id=82311'
$data[8,7] = '//i2.hdslb.com/bfs/openplatform/202404/RFYwkzvt1713951750482.jpeg'
$data[9,0] = '2024-06-22'
$data[9,1] = '合肥·城市动漫节'
$data[9,2] = '包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心'
$data[9,3] = '2024.06.22 10:00-06.23 16:30'
$data[9,4] = 586
$data[9,5] = 50
$data[9,6] = 'https://show.bilibili.com/platform/detail.html?id=85000'
$data[9,7] = '//i2.hdslb.com/bfs/openplatform/202404/U2EZscfQ1714448575403.jpeg'
$data[10,0] = '2024-06-30'
$data[10,1] = '安徽·THO4·隙间皖韵之梦'
$data[10,2] = '北二环与新蚌埠路交汇处 蓝金湾大酒店'
$data[10,3] = '2024.06.30 10:00-06.30 17:00'
$data[10,4] = 18
$data[10,5] = 65
$data[10,6] = 'https://show.bilibili.com/platform/detail.html?id=85119'
$data[10,7] = '//i2.hdslb.com/bfs/openplatform/202405/kuuarwvJ1714932457216.jpeg'
$data[11,0] = '2024-07-20'
$data[11,1] = '合肥·W·A首届童年怀旧only'
$data[11,2] = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$data[11,3] = '2024.07.20 09:30-07.20 17:00'
$data[11,4] = 166
$data[11,5] = 78
$data[11,6] = 'https://show.bilibili.com/platform/detail.html?id=84794'
$data[11,7] = '//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png'
$data[12,0] = '2024-07-20'
$data[12,1] = '安徽·赛马娘Only 2.0'
$data[12,2] = '文忠路1865号 赫拉诺言艺术中心'
$data[12,3] = '2024.07.20 09:00-07.20 17:00'
$data[12,4] = 1
$data[12,5] = 78
$data[12,6] = 'https://show.bilibili.com/platform/detail.html?id=84539'
$data[12,7] = '//i2.hdslb.com/bfs/openplatform/202405/oa09dctb1715068234778.png'
$data[13,0] = '2024-07-27'
$data[13,1] = '安徽·MAX特摄only展'
$data[13,2] = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$data[13,3] = '2024.07.27 09:30-07.27 18:00'
$data[13,4] = 196
$data[13,5] = 50
$data[13,6] = 'https://show.bilibili.com/platform/detail.html?id=83684'
$data[13,7] = '//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg'
$data[14,0] = '2024-08-03'
$data[14,1] = '合肥·第七届环形宇宙动漫游戏嘉年华'
$data[14,2] = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$data[14,3] = '2024.08.03 09:30-08.04 17:00'
$data[14,4] = 710
$data[14,5] = 49
$data[14,6] = 'https://show.bilibili.com/platform/detail.html?id=84767'
$data[14,7] = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'
$ws1.Range("B2:I16").Value = $data

# column A (#) sequence numbers
$colA = New-Object 'object[,]' 15,1
$colA[0,0] = 1
$colA[1,0] = 2
$colA[2,0] = 3
$colA[3,0] = 4
$colA[4,0] = 5
$colA[5,0] = 6
$colA[6,0] = 7
$colA[7,0] = 8
$colA[8,0] = 9
$colA[9,0] = 10
$colA[10,0] = 11
$colA[11,0] = 12
$colA[12,0] = 13
$colA[13,0] = 14
$colA[14,0] = 15
$ws1.Range("A2:A16").Value = $colA

# copy style/format of column A from an existing row onto newly added rows
$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("A15:A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2 = $wb.Worksheets.Item("演出")
# ---- populate 2 data rows (rows 2..3) ----
$ws2.Range("B2:B3").NumberFormat = "@"
$ws2.Range("E2:E3").NumberFormat = "@"
$data = New-Object 'object[,]' 2,8
$data[0,0] = '2024-06-01'
$data[0,1] = '合肥·跨越二次元ACG神级动漫世界巡回演唱会'
$data[0,2] = '金寨路310号合柴1972文创园区C-1号 合肥九莱福'
$data[0,3] = '2024.06.01 20:00-06.01 21:40'
$data[0,4] = 0
$data[0,5] = '已售罄'
$data[0,6] = 'https://show.bilibili.com/platform/detail.html?id=85179'
$data[0,7] = '//i1.hdslb.com/bfs/openplatform/202405/S1x6SBNF1714972333798.jpeg'
$data[1,0] = '2024-08-03'
$data[1,1] = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
$data[1,2] = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$data[1,3] = '2024.08.03 19:30-08.03 21:00'
$data[1,4] = 22
$data[1,5] = 80
$data[1,6] = 'https://show.bilibili.com/platform/detail.html?id=83556'
$data[1,7] = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'
$ws2.Range("B2:I3").Value = $data

# column A (#) sequence numbers
$colA = New-Object 'object[,]' 2,1
$colA[0,0] = 1
$colA[1,0] = 2
$ws2.Range("A2:A3").Value = $colA

# copy style/format of column A from an existing row onto newly added rows
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A3:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws4 = $wb.Worksheets.Item("全部类型")
# ---- populate 17 data rows (rows 2..18) ----
$ws4.Range("B2:B18").NumberFormat = "@"
$ws4.Range("E2:E18").NumberFormat = "@"
$data = New-Object 'object[,]' 17,8
$data[0,0] = '2024-05-18'
$data[0,1] = '合肥·WA二次元饭局（取消）'
$data[0,2] = '临泉路胜利路交叉路（中环国际大厦对面） 太太满庭芳(胜利路店)'
$data[0,3] = '2024.05.18 14:50-05.18 20:00'
$data[0,4] = 67
$data[0,5] = '不可售'
$data[0,6] = 'https://show.bilibili.com/platform/detail.html?id=83978'
$data[0,7] = '//i1.hdslb.com/bfs/openplatform/202404/wK9Yq9Ta1712657384067.jpeg'
$data[1,0] = '2024-05-18'
$data[1,1] = '合肥·梦时空SPO1动漫展（取消）'
$data[1,2] = '阜阳路16号 银瑞林国际大酒店'
$data[1,3] = '2024.05.18 10:00-05.18 17:00'
$data[1,4] = 126
$data[1,5] = '不可售'
$data[1,6] = 'https://show.bilibili.com/platform/detail.html?id=80207'
$data[1,7] = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'
$data[2,0] = '2024-05-18'
$data[2,1] = '合肥·首届偶活企划——偶像计划-闪耀舞台'
$data[2,2] = '阜阳路16号 银瑞林国际大酒店'
$data[2,3] = '2024.05.18 09:00-05.18 17:00'
$data[2,4] = 75
$data[2,5] = 58
$data[2,6] = 'https://show.bilibili.com/platform/detail.html?id=83891'
$data[2,7] = '//i2.hdslb.com/bfs/openplatform/202404/lfqv8l9Q1712453982625.jpeg'
$data[3,0] = '2024-06-01'
$data[3,1] = '合肥·TH元气动漫游戏嘉年华'
$data[3,2] = '北一环路与胜利路交口西北侧中星城2号楼(地铁1号线长淮站D出口） 曙光薇酒店(合肥站店)'
$data[3,3] = '2024.06.01 10:00-06.01 17:00'
$data[3,4] = 3
$data[3,5] = 50
$data[3,6] = 'https://show.bilibili.com/platform/detail.html?id=85181'
$data[3,7] = '//i0.hdslb.com/bfs/openplatform/202405/QcP0IlNZ1715064886624.jpeg'
$data[4,0] = '2024-06-01'
$data[4,1] = '合肥·跨越二次元ACG神级动漫世界巡回演唱会'
$data[4,2] = '金寨路310号合柴1972文创园区C-1号 合肥九莱福'
$data[4,3] = '2024.06.01 20:00-06.01 21:40'
$data[4,4] = 0
$data[4,5] = '已售罄'
$data[4,6] = 'https://show.bilibili.com/platform/detail.html?id=85179'
$data[4,7] = '//i1.hdslb.com/bfs/openplatform/202405/S1x6SBNF1714972333798.jpeg'
$data[5,0] = '2024-06-01'
$data[5,1] = '合肥·运动番only·群青日和'
$data[5,2] = '金寨路287号 合肥明星运动公园'
$data[5,3] = '2024.06.01 09:30-06.01 17:30'
$data[5,4] = 539
$data[5,5] = 80
$data[5,6] = 'https://show.bilibili.com/platform/detail.html?id=83058'
$data[5,7] = '//i2.hdslb.com/bfs/openplatform/202404/Jzeq47lD1714026878824.jpeg'
$data[6,0] = '2024-06-08'
$data[6,1] = '合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~'
$data[6,2] = '锦绣大道3899号 合肥滨湖会展中心'
$data[6,3] = '2024.06.08 09:30-06.09 17:00'
$data[6,4] = 7523
$data[6,5] = 65
$data[6,6] = 'https://show.bilibili.com/platform/detail.html?id=83518'
$data[6,7] = '//i1.hdslb.com/bfs/openplatform/202403/1Sqp42gM1711691520194.jpeg'
$data[7,0] = '2024-06-09'
$data[7,1] = '合肥·第二届华盟动漫次元嘉年华'
$data[7,2] = '常青街道十五里河村合柴1972院内 合肥当代美术馆'
$data[7,3] = '2024.06.09 10:00-06.10 17:00'
$data[7,4] = 474
$data[7,5] = 58
$data[7,6] = 'https://show.bilibili.com/platform/detail.html?id=84081'
$data[7,7] = '//i1.hdslb.com/bfs/openplatform/202404/O5LyHE7j1712732240786.jpeg'
$data[8,0] = '2024-06-09'
$data[8,1] = '合肥·第六届环形宇宙动漫游戏嘉年华内场票·赵成晨'
$data[8,2] = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$data[8,3] = '2024.06.09 09:30-06.09 17:00'
$data[8,4] = 194
$data[8,5] = 238
$data[8,6] = 'https://show.bilibili.com/platform/detail.html?id=84863'
$data[8,7] = '//i1.hdslb.com/bfs/openplatform/202404/I5S4Ih2M1714031127805.jpeg'
$data[9,0] = '2024-06-22'
$data[9,1] = '合肥·Look Look动漫嘉年华'
$data[9,2] = '新站区东方大道288号 少荃体育中心'
$data[9,3] = '2024.06.22 10:00-06.22 17:30'
$data[9,4] = 1069
$data[9,5] = 58
$data[9,6] = 'https://show.bilibili.com/platform/detail.html?id=82311'
$data[9,7] = '//i2.hdslb.com/bfs/openplatform/202404/RFYwkzvt1713951750482.jpeg'
$data[10,0] = '2024-06-22'
$data[10,1] = '合肥·城市动漫节'
$data[10,2] = '包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心'
$data[10,3] = '2024.06.22 10:00-06.23 16:30'
$data[10,4] = 586
$data[10,5] = 50
$data[10,6] = 'https://show.bilibili.com/platform/detail.html?id=85000'
$data[10,7] = '//i2.hdslb.com/bfs/openplatform/202404/U2EZscfQ1714448575403.jpeg'
$data[11,0] = '2024-06-30'
$data[11,1] = '安徽·THO4·隙间皖韵之梦'
$data[11,2] = '北二环与新蚌埠路交汇处 蓝金湾大酒店'
$data[11,3] = '2024.06.30 10:00-06.30 17:00'
$data[11,4] = 18
$data[11,5] = 65
$data[11,6] = 'https://show.bilibili.com/platform/detail.html?id=85119'
$data[11,7] = '//i2.hdslb.com/bfs/openplatform/202405/kuuarwvJ1714932457216.jpeg'
$data[12,0] = '2024-07-20'
$data[12,1] = '合肥·W·A首届童年怀旧only'
$data[12,2] = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$data[12,3] = '2024.07.20 09:30-07.20 17:00'
$data[12,4] = 166
$data[12,5] = 78
$data[12,6] = 'https://show.bilibili.com/platform/detail.html?id=84794'
$data[12,7] = '//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png'
$data[13,0] = '2024-07-20'
$data[13,1] = '安徽·赛马娘Only 2.0'
$data[13,2] = '文忠路1865号 赫拉诺言艺术中心'
$data[13,3] = '2024.07.20 09:00-07.20 17:00'
$data[13,4] = 1
$data[13,5] = 78
$data[13,6] = 'https://show.bilibili.com/platform/detail.html?id=84539'
$data[13,7] = '//i2.hdslb.com/bfs/openplatform/202405/oa09dctb1715068234778.png'
$data[14,0] = '2024-07-27'
$data[14,1] = '安徽·MAX特摄only展'
$data[14,2] = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$data[14,3] = '2024.07.27 09:30-07.27 18:00'
$data[14,4] = 196
$data[14,5] = 50
$data[14,6] = 'https://show.bilibili.com/platform/detail.html?id=83684'
$data[14,7] = '//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg'
$data[15,0] = '2024-08-03'
$data[15,1] = '合肥·第七届环形宇宙动漫游戏嘉年华'
$data[15,2] = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$data[15,3] = '2024.08.03 09:30-08.04 17:00'
$data[15,4] = 710
$data[15,5] = 49
$data[15,6] = 'https://show.bilibili.com/platform/detail.html?id=84767'
$data[15,7] = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'
$data[16,0] = '2024-08-03'
$data[16,1] = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
$data[16,2] = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$data[16,3] = '2024.08.03 19:30-08.03 21:00'
$data[16,4] = 22
$data[16,5] = 80
$data[16,6] = 'https://show.bilibili.com/platform/detail.html?id=83556'
$data[16,7] = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'
$ws4.Range("B2:I18").Value = $data

# column A (#) sequence numbers
$colA = New-Object 'object[,]' 17,1
$colA[0,0] = 1
$colA[1,0] = 2
$colA[2,0] = 3
$colA[3,0] = 4
$colA[4,0] = 5
$colA[5,0] = 6
$colA[6,0] = 7
$colA[7,0] = 8
$colA[8,0] = 9
$colA[9,0] = 10
$colA[10,0] = 11
$colA[11,0] = 12
$colA[12,0] = 13
$colA[13,0] = 14
$colA[14,0] = 15
$colA[15,0] = 16
$colA[16,0] = 17
$ws4.Range("A2:A18").Value = $colA

# copy style/format of column A from an existing row onto newly added rows
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("A16:A18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

